# Edit script: adds a "PO Forecast" sheet, renames headers, and fills forecast data.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Requested quantity" headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value2 = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value2 = "Monthly_PO_Qty"

# --- 2. Add the new "PO Forecast" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# --- 3. Write header row, reusing the bold/border style from the existing header cells ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsForecast.Range("A1").Value2 = "ds"
$wsForecast.Range("B1").Value2 = "PO_Forecast"
$wsForecast.Range("C1").Value2 = "yhat_lower"
$wsForecast.Range("D1").Value2 = "yhat_upper"

# --- 4. Apply the date-style (same as column A on the other sheets) to column A of the data rows ---
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A61").PasteSpecial(-4122)

# --- 5. Fill in the forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper) ---
$data = @(
    @(44990.99999999999, 150, -250.6078954694766, 564.1088080191415),
    @(44997.99999999999, 152, -243.812915757228, 565.3927907979044),
    @(45004.99999999999, 154, -219.7109403895931, 584.4762988321504),
    @(45011.99999999999, 155, -270.9818485649732, 559.0312666001884),
    @(45025.99999999999, 159, -247.3351079068102, 569.3229173707226),
    @(45039.99999999999, 162, -235.8151372618638, 578.5737777523658),
    @(45046.99999999999, 164, -270.1784508766081, 577.2847377554599),
    @(45053.99999999999, 166, -244.1271230718218, 582.9355568477537),
    @(45060.99999999999, 167, -226.7176854972072, 584.5358691372577),
    @(45067.99999999999, 169, -204.1895762387808, 560.0143122092137),
    @(45074.99999999999, 171, -241.0955401536781, 569.7713916656533),
    @(45081.99999999999, 172, -229.7038995180715, 605.9118740529345),
    @(45088.99999999999, 174, -245.0504068527566, 592.7206033415679),
    @(45109.99999999999, 179, -221.8444527470265, 582.0209321573318),
    @(45130.99999999999, 184, -218.1429498191668, 596.0869805444859),
    @(45144.99999999999, 188, -218.5060536391856, 610.7449116821531),
    @(45151.99999999999, 190, -204.4466852199109, 567.729578836194),
    @(45158.99999999999, 191, -215.8276474676452, 576.8736925643695),
    @(45165.99999999999, 193, -224.4537692885566, 604.9324942485266),
    @(45186.99999999999, 198, -193.3249617715551, 571.8593222373941),
    @(45193.99999999999, 200, -192.1584948395872, 609.9612619121378),
    @(45200.99999999999, 202, -214.806295792353, 584.3430516941751),
    @(45214.99999999999, 205, -183.6687517706385, 645.0168399834016),
    @(45221.99999999999, 207, -199.8928390991081, 613.882590740065),
    @(45228.99999999999, 208, -203.3712464751236, 596.6500787825605),
    @(45235.99999999999, 210, -184.1571482932552, 615.3743615991089),
    @(45242.99999999999, 212, -190.9550852444267, 625.0137853817711),
    @(45249.99999999999, 214, -197.5487636912597, 638.0778145247619),
    @(45256.99999999999, 215, -187.8618248386521, 617.8300117660276),
    @(45270.99999999999, 219, -143.0525781071484, 628.0726842601082),
    @(45305.99999999999, 227, -170.4306019477474, 665.2185152569886),
    @(45312.99999999999, 229, -183.3071517255414, 667.2424218385895),
    @(45319.99999999999, 231, -200.6985138798912, 613.7021503156516),
    @(45326.99999999999, 232, -180.9648452071107, 650.604289464809),
    @(45333.99999999999, 234, -197.1971400958164, 665.8555631882721),
    @(45347.99999999999, 237, -161.9534603948039, 640.1545848500188),
    @(45361.99999999999, 241, -133.1458565760508, 644.1330788070374),
    @(45375.99999999999, 244, -172.5830423225918, 653.9158798891086),
    @(45417.99999999999, 255, -137.0553342970521, 681.058744520918),
    @(45424.99999999999, 256, -162.5561222230879, 661.2177571897599),
    @(45431.99999999999, 258, -155.2319991759211, 686.7758146314875),
    @(45459.99999999999, 265, -158.8405799910501, 688.6636576130987),
    @(45466.99999999999, 266, -159.2096191185633, 689.2025938242438),
    @(45473.99999999999, 268, -133.1233982442465, 658.7782611931052),
    @(45480.99999999999, 270, -115.7704249167692, 673.170715609945),
    @(45494.99999999999, 273, -141.9110908779619, 684.255151715775),
    @(45543.99999999999, 285, -144.605036115071, 680.1887598585121),
    @(45550.99999999999, 287, -121.8118551303266, 697.6005664176463),
    @(45557.99999999999, 289, -123.5401259242626, 696.6532849912813),
    @(45564.99999999999, 290, -143.4414486643322, 653.4829814105872),
    @(45571.99999999999, 292, -114.1064786554797, 719.1578919376384),
    @(45585.99999999999, 296, -140.4427909382013, 709.2233034831384),
    @(45592.99999999999, 297, -93.37523929648128, 705.1292581744472),
    @(45599.99999999999, 299, -136.2841929493068, 688.2866227362425),
    @(45606.99999999999, 301, -75.850916194337, 725.2527287747047),
    @(45613.99999999999, 302, -101.3877196515476, 718.9166568389173),
    @(45620.99999999999, 304, -69.31330292229246, 717.2380787853631),
    @(45627.99999999999, 306, -109.0741517221118, 722.4774114231493),
    @(45634.99999999999, 308, -85.38687562603502, 730.883719539659),
    @(45641.99999999999, 309, -103.8415247864212, 744.5020467831605)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $wsForecast.Cells.Item($r, 1).Value2 = $row[0]
    $wsForecast.Cells.Item($r, 2).Value2 = $row[1]
    $wsForecast.Cells.Item($r, 3).Value2 = $row[2]
    $wsForecast.Cells.Item($r, 4).Value2 = $row[3]
}

Write-Host "PO Forecast sheet created and headers updated."
